$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Text = " "
$tr1.Text = "First slide"

$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Text = " "
$tr3.Text = "Third slide"
